# guide41_email.xlsx — "Add files via upload"
#
# On sheet "p1" a new blank row is inserted above the (old) row 5, pushing the
# three Q&A rows (old rows 5-7, each: instruction cell in B, "chartn" cell in
# C, screenshot-name cell in D) down by one, and appending one more trailing
# blank row at the bottom so the used range stays A1:D20. The new blank
# spacer row picks up the "title" row's formatting, and the sheet's
# selection is left on B14.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("p1")

# Insert a new blank row above row 5 — shifts the three instruction rows
# (old 5:7 -> new 6:8) down by one and grows the sheet by a row at the
# bottom (old 19 -> new 20), same as Excel's Rows.Insert "shift down".
$ws2.Rows("5:5").Insert() | Out-Null

# The freshly inserted row 4/5 spacer cell (B4) should carry the same look
# as the title cell above it (B3) rather than the plain body style it
# inherited from the row it was copied from.
$ws2.Range("B3").Copy() | Out-Null
$ws2.Range("B4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Leave the selection on B14, matching where the editor ended up.
$ws2.Range("B14").Select() | Out-Null
